# Updates slice-geometry output values (regenerated by upstream simulation)
# that accompanies the figure added to the documentation home page.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2").Value = 0
$ws.Range("Z3").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("Z5").Value = 0
$ws.Range("Z6").Value = 0
$ws.Range("Z7").Value = 0
$ws.Range("Z8").Value = 0
$ws.Range("Z9").Value = 0
$ws.Range("Z10").Value = 0
$ws.Range("Z11").Value = 0
$ws.Range("Z12").Value = 0
$ws.Range("Z13").Value = 0
$ws.Range("Z14").Value = 0
$ws.Range("Z15").Value = 0
$ws.Range("Z16").Value = 0
$ws.Range("Z17").Value = 0
$ws.Range("E18").Value = 229.8643850819121
$ws.Range("F18").Value = 35.15875917533425
$ws.Range("H18").Value = 224.9321925409561
$ws.Range("I18").Value = 30.49667829853728
$ws.Range("K18").Value = 57.37328724188009
$ws.Range("L18").Value = 9.864385081912133
$ws.Range("M18").Value = 41.76689339905241
$ws.Range("N18").Value = 13.22550038612207
$ws.Range("Q18").Value = 9.503321701462724
$ws.Range("R18").Value = 66431.09432518722
$ws.Range("V18").Value = 224.9321925409561
$ws.Range("Z18").Value = 0
$ws.Range("AF18").Value = 77.37288047629373
$ws.Range("AG18").Value = 46.87620217775645
$ws.Range("AH18").Value = 2925.075015892003
$ws.Range("B19").Value = 229.8643850819121
$ws.Range("C19").Value = 35.15875917533425
$ws.Range("E19").Value = 239.7287701638243
$ws.Range("F19").Value = 46.33548952684643
$ws.Range("H19").Value = 234.7965776228682
$ws.Range("I19").Value = 40.39739547464413
$ws.Range("K19").Value = 62.63318313075887
$ws.Range("L19").Value = 9.864385081912161
$ws.Range("M19").Value = 48.44374871983539
$ws.Range("N19").Value = 14.87044587012606
$ws.Range("P19").Value = 23.60260452535587
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = 53586.42280987477
$ws.Range("V19").Value = 234.7965776228682
$ws.Range("Z19").Value = 0
$ws.Range("AF19").Value = 78.17269548293525
$ws.Range("AG19").Value = 37.77530000829113
$ws.Range("AH19").Value = 2357.178720517366
$ws.Range("B20").Value = 239.7287701638243
$ws.Range("C20").Value = 46.33548952684643
$ws.Range("E20").Value = 249.5931552457365
$ws.Range("F20").Value = 61.175924352853
$ws.Range("H20").Value = 244.6609627047804
$ws.Range("I20").Value = 53.15920023850424
$ws.Range("K20").Value = 68.85750549896575
$ws.Range("L20").Value = 9.864385081912161
$ws.Range("M20").Value = 56.15096057258316
$ws.Range("N20").Value = 17.7096360831204
$ws.Range("P20").Value = 10.84079976149576
$ws.Range("R20").Value = 38479.9400261671
$ws.Range("V20").Value = 244.6609627047804
$ws.Range("Z20").Value = 0
$ws.Range("AF20").Value = 78.97251048957679
$ws.Range("AG20").Value = 25.81331025107255
$ws.Range("AH20").Value = 1610.750559666927
$ws.Range("B21").Value = 249.5931552457365
$ws.Range("C21").Value = 61.175924352853
$ws.Range("E21").Value = 259.4575403276486
$ws.Range("F21").Value = 83.95181749179723
$ws.Range("H21").Value = 254.5253477866925
$ws.Range("I21").Value = 70.96737624597222
$ws.Range("K21").Value = 77.48368812298611
$ws.Range("L21").Value = 9.864385081912133
$ws.Range("M21").Value = 65.8827570725098
$ws.Range("N21").Value = 24.14160699618116
$ws.Range("O21").Value = 13.03262375402778
$ws.Range("R21").Value = 16712.6465138627
$ws.Range("V21").Value = 254.5253477866925
$ws.Range("Y21").Value = 0
$ws.Range("Z21").Value = 0
$ws.Range("AF21").Value = 79.12700463715591
$ws.Range("AG21").Value = 8.159628391183688
$ws.Range("AH21").Value = 509.1608116098621
